# Insert a brand-new weekly price record at row 707 ("Hortaliza, Femacal de
# La Calera - Cebollín" sheet). All existing rows from 707 downward shift
# down by one (707->708, ..., 803->804), and the new row carries the latest
# observation (Fecha serial 45154 = 2023-08-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 707:803 down to 708:804 by inserting a new row at 707.
$ws.Range("A707:R707").EntireRow.Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A707").Value = 3
$ws.Range("B707").Value = "Femacal de La Calera"
$ws.Range("C707").Value = "Coquimbo"
$ws.Range("D707").Value = 45154
$ws.Range("E707").Value = 5
$ws.Range("F707").Value = 100112037
$ws.Range("G707").Value = "Cebollín"
$ws.Range("H707").Value = "Sin especificar"
$ws.Range("I707").Value = "Primera"
$ws.Range("J707").Value = 200
$ws.Range("K707").Value = 3800
$ws.Range("L707").Value = 4000
$ws.Range("M707").Value = 3910
$ws.Range("N707").Value = "`$/paquete 36 unidades"
$ws.Range("O707").Value = "Provincia de Quillota"
$ws.Range("P707").Value = 109
$ws.Range("Q707").Value = 36
$ws.Range("R707").Value = "Hortaliza"
